$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shotgun asset (row 2) dimension swapped: width/depth figures reordered
# "70cm x 16cm x 22cm" -> "70cm x 22cm x 16cm"
$ws.Range("F2").Value = "70cm x 22cm x 16cm"

# Leave the cursor/selection where the author left off saving (cell F2)
$ws.Range("F2").Select()
